$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new resale-data row for 2025-02-01 22:45:37 (row 47)
$row = 47

# Columns A-D hold text that looks numeric/date-like ("2025-02-01", "04").
# Force them to be written as literal text (matching the rest of the sheet)
# by temporarily setting the cell format to Text, then clearing the format
# again afterwards so no extra style index is left on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-02-01"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "22:45:37"

$ws.Cells.Item($row, 3).Value = "Saturday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "04"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 125939
$ws.Cells.Item($row, 6).Value = 141913
$ws.Cells.Item($row, 7).Value = 166763
$ws.Cells.Item($row, 8).Value = 157980
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 142110
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191245
$ws.Cells.Item($row, 14).Value = 115490
$ws.Cells.Item($row, 15).Value = 44846
$ws.Cells.Item($row, 16).Value = 28267
$ws.Cells.Item($row, 17).Value = 63542
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 40787
$ws.Cells.Item($row, 20).Value = -1
